# Merging 20180404_02_002_004 with 20180404_02_005_006, making the trace
# 20180404_02_002_006.
#
# Row 116 held "20180404_02_002_004" and row 117 held "20180404_02_005_006".
# Update row 116's value to the merged trace name, then delete row 117
# entirely so everything below shifts up by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A116").Value = "20180404_02_002_006"
$ws.Rows(117).Delete()

# Leave the view in a state similar to the target (scrolled near the
# merged row, with the row that is now row 117 selected in full).
$ws.Rows(117).Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 115 | Out-Null
